$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.401.13"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.841.32"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.33"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6274"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07446"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2899"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.84"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07710"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").Value = "1.841.57"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6761"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001028"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.82"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "29.387.09"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "233.20"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.31"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.334"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.08%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.23"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.492"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1350"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.37"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07143"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +9.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.462"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.475"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.046"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.037"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.139"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6975"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.573"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.989"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01838"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.813"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").Value = "1.236.05"
$ws.Range("E40").Value = "  -2.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9499"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.47%  "
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "1.991.24"
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.83"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.46"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.728"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.965"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.948"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.35%  "
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3901"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.39%  "
